# Fix Training Data Issue
# Data was off by one day due to how NBA stats were shown; correct the date
# column and recompute the dependent stat/rank values for that corrected date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Date column keeps being stored as text (not auto-converted to a
# serial date number) once we write an ISO-formatted date string into it.
$ws.Range("BF2:BF31").NumberFormat = "@"

$numericCells = @{}
$numericCells["D2"] = 21
$numericCells["F2"] = 10
$numericCells["G2"] = 0.524
$numericCells["I2"] = 37.8
$numericCells["J2"] = 81.59999999999999
$numericCells["K2"] = 0.463
$numericCells["L2"] = 8
$numericCells["M2"] = 22.2
$numericCells["N2"] = 0.362
$numericCells["O2"] = 16
$numericCells["P2"] = 21.6
$numericCells["Q2"] = 0.742
$numericCells["R2"] = 8.4
$numericCells["S2"] = 32
$numericCells["T2"] = 40.5
$numericCells["W2"] = 8.300000000000001
$numericCells["Y2"] = 4.2
$numericCells["Z2"] = 18.6
$numericCells["AB2"] = 99.59999999999999
$numericCells["AC2"] = 0.5
$numericCells["AD2"] = 9
$numericCells["AE2"] = 11
$numericCells["AF2"] = 11
$numericCells["AI2"] = 13
$numericCells["AJ2"] = 21
$numericCells["AK2"] = 7
$numericCells["AN2"] = 13
$numericCells["AO2"] = 22
$numericCells["AS2"] = 14
$numericCells["AV2"] = 10
$numericCells["AW2"] = 9
$numericCells["AX2"] = 20
$numericCells["AY2"] = 9
$numericCells["D3"] = 22
$numericCells["F3"] = 12
$numericCells["G3"] = 0.455
$numericCells["K3"] = 0.454
$numericCells["L3"] = 6.5
$numericCells["M3"] = 18.3
$numericCells["N3"] = 0.353
$numericCells["P3"] = 20.2
$numericCells["Q3"] = 0.766
$numericCells["R3"] = 10.4
$numericCells["S3"] = 31
$numericCells["T3"] = 41.4
$numericCells["U3"] = 18.4
$numericCells["V3"] = 16.4
$numericCells["W3"] = 7.2
$numericCells["X3"] = 4.8
$numericCells["Y3"] = 4.8
$numericCells["AA3"] = 18.8
$numericCells["AB3"] = 95.2
$numericCells["AC3"] = -0.5
$numericCells["AD3"] = 2
$numericCells["AG3"] = 18
$numericCells["AM3"] = 25
$numericCells["AN3"] = 16
$numericCells["AP3"] = 26
$numericCells["AS3"] = 19
$numericCells["AT3"] = 22
$numericCells["AU3"] = 28
$numericCells["AV3"] = 24
$numericCells["AX3"] = 14
$numericCells["AY3"] = 15
$numericCells["D4"] = 20
$numericCells["E4"] = 6
$numericCells["G4"] = 0.3
$numericCells["I4"] = 34.2
$numericCells["J4"] = 79.3
$numericCells["K4"] = 0.431
$numericCells["L4"] = 6.3
$numericCells["M4"] = 18.4
$numericCells["N4"] = 0.341
$numericCells["O4"] = 19.8
$numericCells["Q4"] = 0.75
$numericCells["R4"] = 10.8
$numericCells["S4"] = 30.4
$numericCells["W4"] = 6.9
$numericCells["X4"] = 4.4
$numericCells["Y4"] = 4.4
$numericCells["Z4"] = 22.5
$numericCells["AB4"] = 94.40000000000001
$numericCells["AC4"] = -7.9
$numericCells["AD4"] = 16
$numericCells["AE4"] = 25
$numericCells["AG4"] = 26
$numericCells["AH4"] = 12
$numericCells["AK4"] = 22
$numericCells["AM4"] = 23
$numericCells["AO4"] = 5
$numericCells["AT4"] = 23
$numericCells["AW4"] = 26
$numericCells["AX4"] = 25
$numericCells["AY4"] = 11
$numericCells["AZ4"] = 25
$numericCells["BA4"] = 9
$numericCells["BC4"] = 28
$numericCells["AD5"] = 9
$numericCells["AF5"] = 17
$numericCells["AG5"] = 15
$numericCells["AT5"] = 14
$numericCells["AU5"] = 27
$numericCells["AX5"] = 8
$numericCells["AZ5"] = 7
$numericCells["D6"] = 18
$numericCells["F6"] = 10
$numericCells["G6"] = 0.444
$numericCells["I6"] = 35.4
$numericCells["K6"] = 0.431
$numericCells["L6"] = 5.3
$numericCells["N6"] = 0.336
$numericCells["O6"] = 18
$numericCells["P6"] = 23.1
$numericCells["Q6"] = 0.781
$numericCells["S6"] = 33.8
$numericCells["T6"] = 46.4
$numericCells["U6"] = 22.1
$numericCells["V6"] = 16.3
$numericCells["X6"] = 4.3
$numericCells["Y6"] = 6.5
$numericCells["Z6"] = 19.7
$numericCells["AA6"] = 21.8
$numericCells["AB6"] = 94.2
$numericCells["AC6"] = 0.2
$numericCells["AD6"] = 30
$numericCells["AF6"] = 11
$numericCells["AI6"] = 24
$numericCells["AK6"] = 23
$numericCells["AN6"] = 22
$numericCells["AO6"] = 13
$numericCells["AP6"] = 14
$numericCells["AQ6"] = 7
$numericCells["AT6"] = 3
$numericCells["AV6"] = 23
$numericCells["AX6"] = 28
$numericCells["AY6"] = 28
$numericCells["AZ6"] = 9
$numericCells["BA6"] = 8
$numericCells["BB6"] = 25
$numericCells["D7"] = 20
$numericCells["E7"] = 7
$numericCells["G7"] = 0.35
$numericCells["H7"] = 48.8
$numericCells["I7"] = 35
$numericCells["J7"] = 84.40000000000001
$numericCells["K7"] = 0.415
$numericCells["L7"] = 6.6
$numericCells["M7"] = 19.2
$numericCells["N7"] = 0.342
$numericCells["Q7"] = 0.744
$numericCells["R7"] = 12.1
$numericCells["S7"] = 31.9
$numericCells["T7"] = 44
$numericCells["U7"] = 18.4
$numericCells["V7"] = 16.1
$numericCells["W7"] = 7.9
$numericCells["X7"] = 4.7
$numericCells["Y7"] = 5.9
$numericCells["Z7"] = 19.7
$numericCells["AA7"] = 19.3
$numericCells["AB7"] = 92.40000000000001
$numericCells["AC7"] = -7.3
$numericCells["AD7"] = 16
$numericCells["AE7"] = 22
$numericCells["AF7"] = 23
$numericCells["AG7"] = 23
$numericCells["AI7"] = 26
$numericCells["AJ7"] = 9
$numericCells["AK7"] = 29
$numericCells["AQ7"] = 18
$numericCells["AT7"] = 10
$numericCells["AU7"] = 29
$numericCells["AW7"] = 16
$numericCells["AX7"] = 17
$numericCells["AY7"] = 24
$numericCells["AZ7"] = 8
$numericCells["BB7"] = 27
$numericCells["AD8"] = 2
$numericCells["AK8"] = 9
$numericCells["AO8"] = 16
$numericCells["AV8"] = 11
$numericCells["AX8"] = 24
$numericCells["AY8"] = 4
$numericCells["AZ8"] = 16
$numericCells["BB8"] = 7
$numericCells["AD9"] = 9
$numericCells["AK9"] = 13
$numericCells["AL9"] = 13
$numericCells["AM9"] = 17
$numericCells["AN9"] = 10
$numericCells["AS9"] = 9
$numericCells["AW9"] = 23
$numericCells["AZ9"] = 20
$numericCells["D10"] = 21
$numericCells["F10"] = 11
$numericCells["G10"] = 0.476
$numericCells["I10"] = 37.9
$numericCells["K10"] = 0.45
$numericCells["L10"] = 6.3
$numericCells["M10"] = 19.8
$numericCells["N10"] = 0.32
$numericCells["O10"] = 17.2
$numericCells["P10"] = 25.4
$numericCells["Q10"] = 0.677
$numericCells["R10"] = 14.1
$numericCells["S10"] = 30.2
$numericCells["T10"] = 44.4
$numericCells["U10"] = 20.1
$numericCells["V10"] = 15.8
$numericCells["W10"] = 9.800000000000001
$numericCells["Y10"] = 5
$numericCells["Z10"] = 20.4
$numericCells["AA10"] = 20.5
$numericCells["AB10"] = 99.3
$numericCells["AC10"] = 0
$numericCells["AD10"] = 9
$numericCells["AF10"] = 17
$numericCells["AG10"] = 15
$numericCells["AI10"] = 11
$numericCells["AJ10"] = 10
$numericCells["AK10"] = 16
$numericCells["AL10"] = 24
$numericCells["AM10"] = 21
$numericCells["AO10"] = 15
$numericCells["AS10"] = 25
$numericCells["AT10"] = 8
$numericCells["AV10"] = 18
$numericCells["AW10"] = 1
$numericCells["AX10"] = 8
$numericCells["BA10"] = 17
$numericCells["BC10"] = 16
$numericCells["AD11"] = 2
$numericCells["AX11"] = 8
$numericCells["AY11"] = 14
$numericCells["BB11"] = 8
$numericCells["BC11"] = 8
$numericCells["AD12"] = 2
$numericCells["AE12"] = 4
$numericCells["AN12"] = 14
$numericCells["AZ12"] = 14
$numericCells["D13"] = 21
$numericCells["E13"] = 18
$numericCells["G13"] = 0.857
$numericCells["I13"] = 36.2
$numericCells["J13"] = 80.2
$numericCells["K13"] = 0.451
$numericCells["L13"] = 7.3
$numericCells["M13"] = 20.2
$numericCells["N13"] = 0.363
$numericCells["O13"] = 18.4
$numericCells["P13"] = 23.2
$numericCells["Q13"] = 0.793
$numericCells["V13"] = 15.7
$numericCells["X13"] = 6.6
$numericCells["Y13"] = 4.6
$numericCells["AA13"] = 22.6
$numericCells["AB13"] = 98.09999999999999
$numericCells["AC13"] = 8.6
$numericCells["AD13"] = 9
$numericCells["AK13"] = 14
$numericCells["AL13"] = 17
$numericCells["AM13"] = 19
$numericCells["AP13"] = 13
$numericCells["AT13"] = 9
$numericCells["AV13"] = 16
$numericCells["AX13"] = 3
$numericCells["AY13"] = 13
$numericCells["AZ13"] = 6
$numericCells["BA13"] = 5
$numericCells["BB13"] = 19
$numericCells["AD14"] = 2
$numericCells["AH14"] = 21
$numericCells["AK14"] = 5
$numericCells["AL14"] = 12
$numericCells["AN14"] = 25
$numericCells["AT14"] = 15
$numericCells["AW14"] = 12
$numericCells["AZ14"] = 27
$numericCells["D15"] = 20
$numericCells["F15"] = 10
$numericCells["G15"] = 0.5
$numericCells["J15"] = 85.90000000000001
$numericCells["K15"] = 0.441
$numericCells["L15"] = 10.6
$numericCells["M15"] = 26.4
$numericCells["N15"] = 0.402
$numericCells["O15"] = 14.7
$numericCells["P15"] = 20.1
$numericCells["Q15"] = 0.729
$numericCells["R15"] = 10.3
$numericCells["S15"] = 33.3
$numericCells["T15"] = 43.6
$numericCells["U15"] = 24
$numericCells["V15"] = 15
$numericCells["X15"] = 6.1
$numericCells["Y15"] = 4.4
$numericCells["Z15"] = 21.6
$numericCells["AA15"] = 18.3
$numericCells["AB15"] = 101
$numericCells["AC15"] = -2
$numericCells["AD15"] = 16
$numericCells["AF15"] = 11
$numericCells["AG15"] = 13
$numericCells["AI15"] = 12
$numericCells["AJ15"] = 4
$numericCells["AN15"] = 3
$numericCells["AO15"] = 28
$numericCells["AP15"] = 27
$numericCells["AS15"] = 8
$numericCells["AT15"] = 13
$numericCells["AV15"] = 15
$numericCells["AY15"] = 11
$numericCells["AZ15"] = 21
$numericCells["BA15"] = 29
$numericCells["BC15"] = 21
$numericCells["AD16"] = 16
$numericCells["AG16"] = 13
$numericCells["AH16"] = 12
$numericCells["AK16"] = 15
$numericCells["AN16"] = 23
$numericCells["AP16"] = 25
$numericCells["AR16"] = 18
$numericCells["AS16"] = 21
$numericCells["AW16"] = 24
$numericCells["AZ16"] = 10
$numericCells["BA16"] = 20
$numericCells["BB16"] = 26
$numericCells["D17"] = 21
$numericCells["F17"] = 5
$numericCells["G17"] = 0.762
$numericCells["I17"] = 38.2
$numericCells["J17"] = 75.2
$numericCells["K17"] = 0.508
$numericCells["L17"] = 8.4
$numericCells["N17"] = 0.395
$numericCells["O17"] = 18.7
$numericCells["P17"] = 24.5
$numericCells["Q17"] = 0.765
$numericCells["S17"] = 29.2
$numericCells["T17"] = 35.6
$numericCells["U17"] = 23.9
$numericCells["V17"] = 15.9
$numericCells["W17"] = 9.300000000000001
$numericCells["Y17"] = 2.9
$numericCells["AB17"] = 103.5
$numericCells["AC17"] = 7.7
$numericCells["AD17"] = 9
$numericCells["AL17"] = 9
$numericCells["AN17"] = 5
$numericCells["AV17"] = 19
$numericCells["AW17"] = 4
$numericCells["BB17"] = 6
$numericCells["D18"] = 20
$numericCells["E18"] = 4
$numericCells["G18"] = 0.2
$numericCells["H18"] = 48.8
$numericCells["I18"] = 34.2
$numericCells["J18"] = 81.7
$numericCells["K18"] = 0.419
$numericCells["L18"] = 7.5
$numericCells["M18"] = 20
$numericCells["N18"] = 0.376
$numericCells["O18"] = 14.9
$numericCells["Q18"] = 0.77
$numericCells["S18"] = 29.5
$numericCells["U18"] = 21.4
$numericCells["X18"] = 5.2
$numericCells["Y18"] = 5.3
$numericCells["Z18"] = 21.1
$numericCells["AA18"] = 19.7
$numericCells["AB18"] = 90.8
$numericCells["AC18"] = -8.6
$numericCells["AD18"] = 16
$numericCells["AE18"] = 29
$numericCells["AI18"] = 28
$numericCells["AJ18"] = 20
$numericCells["AK18"] = 28
$numericCells["AL18"] = 16
$numericCells["AM18"] = 20
$numericCells["AN18"] = 9
$numericCells["AO18"] = 27
$numericCells["AR18"] = 21
$numericCells["AU18"] = 17
$numericCells["AV18"] = 25
$numericCells["AY18"] = 19
$numericCells["AZ18"] = 17
$numericCells["BA18"] = 22
$numericCells["D19"] = 20
$numericCells["E19"] = 9
$numericCells["G19"] = 0.45
$numericCells["H19"] = 48.3
$numericCells["I19"] = 37.8
$numericCells["J19"] = 89.2
$numericCells["K19"] = 0.423
$numericCells["L19"] = 7.8
$numericCells["N19"] = 0.323
$numericCells["O19"] = 21.5
$numericCells["P19"] = 26.8
$numericCells["Q19"] = 0.8
$numericCells["S19"] = 32.1
$numericCells["U19"] = 22.8
$numericCells["V19"] = 14.8
$numericCells["W19"] = 9.5
$numericCells["X19"] = 3.3
$numericCells["Y19"] = 6.3
$numericCells["Z19"] = 17.3
$numericCells["AA19"] = 22.7
$numericCells["AB19"] = 104.7
$numericCells["AC19"] = 2.8
$numericCells["AD19"] = 16
$numericCells["AE19"] = 18
$numericCells["AF19"] = 17
$numericCells["AG19"] = 19
$numericCells["AI19"] = 14
$numericCells["AJ19"] = 1
$numericCells["AL19"] = 14
$numericCells["AN19"] = 27
$numericCells["AS19"] = 13
$numericCells["AU19"] = 8
$numericCells["AV19"] = 14
$numericCells["AW19"] = 2
$numericCells["BA19"] = 4
$numericCells["BB19"] = 4
$numericCells["BC19"] = 9
$numericCells["AD20"] = 24
$numericCells["AE20"] = 18
$numericCells["AL20"] = 23
$numericCells["AN20"] = 6
$numericCells["AP20"] = 12
$numericCells["AS20"] = 26
$numericCells["AW20"] = 3
$numericCells["AY20"] = 29
$numericCells["AZ20"] = 24
$numericCells["BB20"] = 8
$numericCells["D21"] = 19
$numericCells["F21"] = 14
$numericCells["G21"] = 0.263
$numericCells["I21"] = 35.8
$numericCells["J21"] = 83.3
$numericCells["K21"] = 0.43
$numericCells["L21"] = 8.9
$numericCells["M21"] = 25.3
$numericCells["N21"] = 0.351
$numericCells["O21"] = 13.7
$numericCells["P21"] = 17.8
$numericCells["R21"] = 10.7
$numericCells["S21"] = 28.7
$numericCells["T21"] = 39.4
$numericCells["U21"] = 19.7
$numericCells["V21"] = 12.8
$numericCells["W21"] = 8.199999999999999
$numericCells["X21"] = 4.8
$numericCells["Y21"] = 4.2
$numericCells["Z21"] = 23.2
$numericCells["AA21"] = 18.7
$numericCells["AC21"] = -3.8
$numericCells["AD21"] = 24
$numericCells["AF21"] = 25
$numericCells["AK21"] = 24
$numericCells["AS21"] = 29
$numericCells["AW21"] = 11
$numericCells["AX21"] = 13
$numericCells["AY21"] = 7
$numericCells["AZ21"] = 29
$numericCells["BA21"] = 28
$numericCells["D22"] = 19
$numericCells["E22"] = 15
$numericCells["G22"] = 0.789
$numericCells["J22"] = 82.5
$numericCells["K22"] = 0.464
$numericCells["L22"] = 6.2
$numericCells["N22"] = 0.336
$numericCells["P22"] = 26.8
$numericCells["Q22"] = 0.819
$numericCells["S22"] = 35.5
$numericCells["T22"] = 46.7
$numericCells["U22"] = 20.9
$numericCells["Y22"] = 4.3
$numericCells["Z22"] = 21.9
$numericCells["AA22"] = 20.6
$numericCells["AB22"] = 104.7
$numericCells["AC22"] = 5.6
$numericCells["AD22"] = 24
$numericCells["AE22"] = 4
$numericCells["AH22"] = 10
$numericCells["AK22"] = 6
$numericCells["AN22"] = 20
$numericCells["AQ22"] = 1
$numericCells["AU22"] = 18
$numericCells["AV22"] = 22
$numericCells["AW22"] = 17
$numericCells["AX22"] = 4
$numericCells["AY22"] = 10
$numericCells["BB22"] = 3
$numericCells["AD23"] = 9
$numericCells["AE23"] = 25
$numericCells["AF23"] = 27
$numericCells["AH23"] = 8
$numericCells["AJ23"] = 19
$numericCells["AK23"] = 12
$numericCells["AN23"] = 15
$numericCells["AW23"] = 13
$numericCells["AX23"] = 16
$numericCells["AZ23"] = 13
$numericCells["AD24"] = 2
$numericCells["AE24"] = 22
$numericCells["AF24"] = 27
$numericCells["AG24"] = 24
$numericCells["AJ24"] = 2
$numericCells["AO24"] = 21
$numericCells["AP24"] = 15
$numericCells["AT24"] = 4
$numericCells["AZ24"] = 15
$numericCells["BA24"] = 19
$numericCells["BB24"] = 11
$numericCells["BC24"] = 27
$numericCells["D25"] = 20
$numericCells["E25"] = 11
$numericCells["G25"] = 0.55
$numericCells["H25"] = 48.3
$numericCells["I25"] = 37.7
$numericCells["J25"] = 81.5
$numericCells["K25"] = 0.462
$numericCells["L25"] = 9.300000000000001
$numericCells["M25"] = 25.3
$numericCells["N25"] = 0.366
$numericCells["O25"] = 16.6
$numericCells["P25"] = 22.5
$numericCells["Q25"] = 0.738
$numericCells["U25"] = 18.6
$numericCells["V25"] = 16.1
$numericCells["W25"] = 8.300000000000001
$numericCells["X25"] = 4.9
$numericCells["Y25"] = 3.8
$numericCells["AB25"] = 101.2
$numericCells["AC25"] = 1.7
$numericCells["AD25"] = 16
$numericCells["AE25"] = 11
$numericCells["AI25"] = 15
$numericCells["AJ25"] = 22
$numericCells["AK25"] = 8
$numericCells["AO25"] = 17
$numericCells["AQ25"] = 21
$numericCells["AR25"] = 16
$numericCells["AS25"] = 20
$numericCells["AV25"] = 20
$numericCells["AX25"] = 12
$numericCells["AY25"] = 5
$numericCells["BA25"] = 18
$numericCells["AD26"] = 2
$numericCells["AE26"] = 1
$numericCells["AH26"] = 21
$numericCells["AJ26"] = 5
$numericCells["AK26"] = 10
$numericCells["AQ26"] = 2
$numericCells["AU26"] = 7
$numericCells["AD27"] = 24
$numericCells["AE27"] = 25
$numericCells["AF27"] = 23
$numericCells["AG27"] = 25
$numericCells["AK27"] = 21
$numericCells["AL27"] = 18
$numericCells["AN27"] = 21
$numericCells["AQ27"] = 6
$numericCells["AT27"] = 17
$numericCells["BC27"] = 22
$numericCells["D28"] = 19
$numericCells["E28"] = 15
$numericCells["G28"] = 0.789
$numericCells["I28"] = 40.6
$numericCells["K28"] = 0.488
$numericCells["L28"] = 8.300000000000001
$numericCells["M28"] = 20.9
$numericCells["N28"] = 0.397
$numericCells["Q28"] = 0.749
$numericCells["S28"] = 33.9
$numericCells["T28"] = 42.3
$numericCells["U28"] = 24.8
$numericCells["AA28"] = 17.7
$numericCells["AB28"] = 101.8
$numericCells["AC28"] = 8.699999999999999
$numericCells["AD28"] = 24
$numericCells["AE28"] = 4
$numericCells["AL28"] = 10
$numericCells["AM28"] = 18
$numericCells["AN28"] = 4
$numericCells["AV28"] = 12
$numericCells["AW28"] = 14
$numericCells["AX28"] = 18
$numericCells["BB28"] = 12
$numericCells["D29"] = 19
$numericCells["F29"] = 12
$numericCells["G29"] = 0.368
$numericCells["I29"] = 35.4
$numericCells["J29"] = 82.7
$numericCells["K29"] = 0.427
$numericCells["L29"] = 7.1
$numericCells["M29"] = 21.6
$numericCells["N29"] = 0.328
$numericCells["O29"] = 19.9
$numericCells["P29"] = 25.6
$numericCells["Q29"] = 0.778
$numericCells["R29"] = 12.4
$numericCells["S29"] = 30.3
$numericCells["T29"] = 42.6
$numericCells["U29"] = 17.4
$numericCells["V29"] = 14.8
$numericCells["Y29"] = 5.2
$numericCells["Z29"] = 23.1
$numericCells["AA29"] = 22.9
$numericCells["AB29"] = 97.7
$numericCells["AC29"] = -0.7
$numericCells["AD29"] = 24
$numericCells["AE29"] = 22
$numericCells["AG29"] = 22
$numericCells["AH29"] = 4
$numericCells["AI29"] = 25
$numericCells["AK29"] = 26
$numericCells["AL29"] = 19
$numericCells["AN29"] = 24
$numericCells["AO29"] = 4
$numericCells["AQ29"] = 8
$numericCells["AS29"] = 24
$numericCells["AT29"] = 16
$numericCells["AV29"] = 13
$numericCells["AW29"] = 22
$numericCells["AX29"] = 19
$numericCells["AZ29"] = 28
$numericCells["BA29"] = 3
$numericCells["BB29"] = 20
$numericCells["BC29"] = 18
$numericCells["AE30"] = 29
$numericCells["AK30"] = 25
$numericCells["AN30"] = 26
$numericCells["AO30"] = 17
$numericCells["AP30"] = 16
$numericCells["AQ30"] = 20
$numericCells["AS30"] = 30
$numericCells["BA30"] = 10
$numericCells["AD31"] = 16
$numericCells["AE31"] = 18
$numericCells["AF31"] = 17
$numericCells["AP31"] = 22
$numericCells["AR31"] = 20
$numericCells["AV31"] = 17
$numericCells["AX31"] = 15
$numericCells["AY31"] = 6
$numericCells["BA31"] = 20
$numericCells["BC31"] = 19

foreach ($key in $numericCells.Keys) {
    $ws.Range($key).Value = $numericCells[$key]
}

# Corrected date text (was off by one day): "12-10-2013-14" -> "2013-12-10"
$dateCells = @{}
$dateCells["BF2"] = "2013-12-10"
$dateCells["BF3"] = "2013-12-10"
$dateCells["BF4"] = "2013-12-10"
$dateCells["BF5"] = "2013-12-10"
$dateCells["BF6"] = "2013-12-10"
$dateCells["BF7"] = "2013-12-10"
$dateCells["BF8"] = "2013-12-10"
$dateCells["BF9"] = "2013-12-10"
$dateCells["BF10"] = "2013-12-10"
$dateCells["BF11"] = "2013-12-10"
$dateCells["BF12"] = "2013-12-10"
$dateCells["BF13"] = "2013-12-10"
$dateCells["BF14"] = "2013-12-10"
$dateCells["BF15"] = "2013-12-10"
$dateCells["BF16"] = "2013-12-10"
$dateCells["BF17"] = "2013-12-10"
$dateCells["BF18"] = "2013-12-10"
$dateCells["BF19"] = "2013-12-10"
$dateCells["BF20"] = "2013-12-10"
$dateCells["BF21"] = "2013-12-10"
$dateCells["BF22"] = "2013-12-10"
$dateCells["BF23"] = "2013-12-10"
$dateCells["BF24"] = "2013-12-10"
$dateCells["BF25"] = "2013-12-10"
$dateCells["BF26"] = "2013-12-10"
$dateCells["BF27"] = "2013-12-10"
$dateCells["BF28"] = "2013-12-10"
$dateCells["BF29"] = "2013-12-10"
$dateCells["BF30"] = "2013-12-10"
$dateCells["BF31"] = "2013-12-10"

foreach ($key in $dateCells.Keys) {
    $ws.Range($key).Value = $dateCells[$key]
}
